$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1806"
$ws.Range("E17").Value = "1807"
$ws.Range("E18").Value = "1809"
$ws.Range("E19").Value = "1810"
$ws.Range("E20").Value = "1811"
$ws.Range("E21").Value = "1812"
